$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny floating-point nudge on the previous last row's timestamp value.
$ws.Cells.Item(94, 1).Value = 44407.76788895139

# Append the new row of data.
$ws.Cells.Item(95, 1).Value = 44408.76765188529
$ws.Cells.Item(95, 2).Value = 80421
$ws.Cells.Item(95, 3).Value = 67909
$ws.Cells.Item(95, 4).Value = 3653
$ws.Cells.Item(95, 5).Value = 2271
$ws.Cells.Item(95, 6).Value = 1648
$ws.Cells.Item(95, 7).Value = 21223
$ws.Cells.Item(95, 8).Value = 1651
$ws.Cells.Item(95, 9).Value = 917
$ws.Cells.Item(95, 10).Value = 201
